$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the "format" sheet (with all data/styles) to the end of the workbook.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2.Name = "Sheet1"

# Update the header/definition cells on the new sheet to the "obj." prefixed
# definitions (bug fix: list defined inside a nested bean).
$ws2.Range("C1").Value = "obj.deflist#key?listToPropKey=true"
$ws2.Range("D1").Value = "obj.deflist#value?listToPropValue=true"
$ws2.Range("A3").Value = "obj.deflist#~"

# Make the new sheet the active tab (mirrors the author's saved view state).
$ws2.Activate()
